$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: "Volume 31 Number 7" -> "...8"; week-of dates advance by one week ---
$ws.Range("A8").Value = "Volume 31   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/19/2024  Through  2/25/2024"

# --- Crime Complaints table (rows 14-29): new weekly data refresh ---
# A14 is a stable text-styled (s=14) reference cell used to restore General/text
# formatting on cells that flip from a numeric stat back to the "no activity" placeholder.

# Row 14
$ws.Range("N14").Value = 0
$ws.Range("N14").NumberFormat = "#,##0.0;""-""#,##0.0"
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = -100
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J15").Value = 2
# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 500
$ws.Range("F16").Value = 21
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = 153.846153846154
$ws.Range("L16").Value = 106.25
$ws.Range("M16").Value = -10.810810810810
$ws.Range("N16").Value = -81.967213114754
# Row 17
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 23
$ws.Range("H17").Value = 4.545454545454
$ws.Range("I17").Value = 38
$ws.Range("J17").Value = 43
$ws.Range("K17").Value = -11.627906976744
$ws.Range("L17").Value = 46.153846153846
$ws.Range("M17").Value = 80.952380952380
$ws.Range("N17").Value = -39.682539682539
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 5.882352941176
$ws.Range("I18").Value = 29
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = -17.142857142857
$ws.Range("L18").Value = -27.5
$ws.Range("M18").Value = 314.285714285714
$ws.Range("N18").Value = -74.561403508771
# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 7.547169811320
$ws.Range("I19").Value = 104
$ws.Range("J19").Value = 99
$ws.Range("K19").Value = 5.050505050505
$ws.Range("L19").Value = 15.555555555555
$ws.Range("M19").Value = 82.456140350877
$ws.Range("N19").Value = -32.903225806451
# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 10
$ws.Range("K20").Value = -30
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -93.137254901960
# Row 21
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 125
$ws.Range("G21").Value = 106
$ws.Range("H21").Value = 17.924528301886
$ws.Range("I21").Value = 212
$ws.Range("J21").Value = 203
$ws.Range("K21").Value = 4.433497536945
$ws.Range("L21").Value = 17.127071823204
$ws.Range("M21").Value = 63.076923076923
$ws.Range("N21").Value = -65.806451612903
# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = 66.666666666666
$ws.Range("L22").Value = 0
# Row 23
$ws.Range("C23").Value = 2
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -40
$ws.Range("I23").Value = 7
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = -12.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 133.333333333333
# Row 24
$ws.Range("C24").Value = 42
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = -6.666666666666
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 149
$ws.Range("H24").Value = 25.503355704698
$ws.Range("I24").Value = 300
$ws.Range("J24").Value = 283
$ws.Range("K24").Value = 6.007067137809
$ws.Range("L24").Value = 57.894736842105
$ws.Range("M24").Value = 52.284263959390
# Row 25
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 180
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = 35.714285714285
$ws.Range("I25").Value = 67
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = 11.666666666666
$ws.Range("L25").Value = 45.652173913043
$ws.Range("M25").Value = 4.6875
# Row 26
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = 0
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 1
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("G26").Value = 1
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("H26").Value = 0
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I26").Value = 1
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 2
$ws.Range("K26").Value = -50
$ws.Range("L26").Value = -50
# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = -42.857142857142
$ws.Range("L27").Value = -50
# Row 28
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G28").Value = 2
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = -66.666666666666
$ws.Range("N28").Value = -83.333333333333
# Row 29
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 3
$ws.Range("K29").Value = -66.666666666666
$ws.Range("N29").Value = -80

$excel.CutCopyMode = $false

